$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.824357333333333
$ws.Range("H2").Value = 5.473072
$ws.Range("I2").Value = 0.003038792251822429
$ws.Range("J2").Value = 0.003038792251822428
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.815322333333333
$ws.Range("N2").Value = 8.445967
$ws.Range("O2").Value = 0.4327245671751157
$ws.Range("P2").Value = 0.4327245671751158
$ws.Range("Q2").Value = 5.136153944513778
$ws.Range("R2").Value = 46.225385500624
$ws.Range("S2").Value = 0.001314960061904956
$ws.Range("T2").Value = 0.001314960061904956

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.824357333333333
$ws.Range("H3").Value = 5.473072
$ws.Range("I3").Value = 0.003038792251822429
$ws.Range("J3").Value = 0.003038792251822428
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.808242333333334
$ws.Range("N3").Value = 8.424727
$ws.Range("O3").Value = 0.4316363472227054
$ws.Range("P3").Value = 0.4316363472227054
$ws.Range("Q3").Value = 5.123237494593779
$ws.Range("R3").Value = 46.109137451344
$ws.Range("S3").Value = 0.001311653187545293
$ws.Range("T3").Value = 0.001311653187545293

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.824357333333333
$ws.Range("H4").Value = 5.473072
$ws.Range("I4").Value = 0.003038792251822429
$ws.Range("J4").Value = 0.003038792251822428
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.882473
$ws.Range("N4").Value = 2.647419
$ws.Range("O4").Value = 0.1356390856021789
$ws.Range("P4").Value = 0.1356390856021789
$ws.Range("Q4").Value = 1.609946089018667
$ws.Range("R4").Value = 14.489514801168
$ws.Range("S4").Value = 0.0004121790023721802
$ws.Range("T4").Value = 0.0004121790023721802

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 577.3276773333333
$ws.Range("H5").Value = 1731.983032
$ws.Range("I5").Value = 0.961642130403093
$ws.Range("J5").Value = 0.961642130403093
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.815322333333333
$ws.Range("N5").Value = 8.445967
$ws.Range("O5").Value = 0.4327245671751157
$ws.Range("P5").Value = 0.4327245671751158
$ws.Range("Q5").Value = 1625.363503647994
$ws.Range("R5").Value = 14628.27153283194
$ws.Range("S5").Value = 0.4161261746560346
$ws.Range("T5").Value = 0.4161261746560346

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 577.3276773333333
$ws.Range("H6").Value = 1731.983032
$ws.Range("I6").Value = 0.961642130403093
$ws.Range("J6").Value = 0.961642130403093
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.808242333333334
$ws.Range("N6").Value = 8.424727
$ws.Range("O6").Value = 0.4316363472227054
$ws.Range("P6").Value = 0.4316363472227054
$ws.Range("Q6").Value = 1621.276023692474
$ws.Range("R6").Value = 14591.48421323227
$ws.Range("S6").Value = 0.4150796965026516
$ws.Range("T6").Value = 0.4150796965026516

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 577.3276773333333
$ws.Range("H7").Value = 1731.983032
$ws.Range("I7").Value = 0.961642130403093
$ws.Range("J7").Value = 0.961642130403093
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.882473
$ws.Range("N7").Value = 2.647419
$ws.Range("O7").Value = 0.1356390856021789
$ws.Range("P7").Value = 0.1356390856021789
$ws.Range("Q7").Value = 509.4760873993786
$ws.Range("R7").Value = 4585.284786594408
$ws.Range("S7").Value = 0.1304362592444068
$ws.Range("T7").Value = 0.1304362592444068

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.8008510000000001
$ws.Range("H8").Value = 2.402553
$ws.Range("I8").Value = 0.001333960057713973
$ws.Range("J8").Value = 0.001333960057713973
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.815322333333333
$ws.Range("N8").Value = 8.445967
$ws.Range("O8").Value = 0.4327245671751157
$ws.Range("P8").Value = 0.4327245671751158
$ws.Range("Q8").Value = 2.254653705972333
$ws.Range("R8").Value = 20.291883353751
$ws.Range("S8").Value = 0.0005772372886031715
$ws.Range("T8").Value = 0.0005772372886031715

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.8008510000000001
$ws.Range("H9").Value = 2.402553
$ws.Range("I9").Value = 0.001333960057713973
$ws.Range("J9").Value = 0.001333960057713973
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.808242333333334
$ws.Range("N9").Value = 8.424727
$ws.Range("O9").Value = 0.4316363472227054
$ws.Range("P9").Value = 0.4316363472227054
$ws.Range("Q9").Value = 2.248983680892334
$ws.Range("R9").Value = 20.240853128031
$ws.Range("S9").Value = 0.0005757856466526487
$ws.Range("T9").Value = 0.0005757856466526487

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8008510000000001
$ws.Range("H10").Value = 2.402553
$ws.Range("I10").Value = 0.001333960057713973
$ws.Range("J10").Value = 0.001333960057713973
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.882473
$ws.Range("N10").Value = 2.647419
$ws.Range("O10").Value = 0.1356390856021789
$ws.Range("P10").Value = 0.1356390856021789
$ws.Range("Q10").Value = 0.706729384523
$ws.Range("R10").Value = 6.360564460707
$ws.Range("S10").Value = 0.0001809371224581531
$ws.Range("T10").Value = 0.0001809371224581531

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 20.403171
$ws.Range("H11").Value = 61.209513
$ws.Range("I11").Value = 0.03398511728737064
$ws.Range("J11").Value = 0.03398511728737064
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.815322333333333
$ws.Range("N11").Value = 8.445967
$ws.Range("O11").Value = 0.4327245671751157
$ws.Range("P11").Value = 0.4327245671751158
$ws.Range("Q11").Value = 57.44150298711899
$ws.Range("R11").Value = 516.973526884071
$ws.Range("S11").Value = 0.014706195168573
$ws.Range("T11").Value = 0.014706195168573

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 20.403171
$ws.Range("H12").Value = 61.209513
$ws.Range("I12").Value = 0.03398511728737064
$ws.Range("J12").Value = 0.03398511728737064
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2.808242333333334
$ws.Range("N12").Value = 8.424727
$ws.Range("O12").Value = 0.4316363472227054
$ws.Range("P12").Value = 0.4316363472227054
$ws.Range("Q12").Value = 57.29704853643901
$ws.Range("R12").Value = 515.673436827951
$ws.Range("S12").Value = 0.01466921188585588
$ws.Range("T12").Value = 0.01466921188585588

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 20.403171
$ws.Range("H13").Value = 61.209513
$ws.Range("I13").Value = 0.03398511728737064
$ws.Range("J13").Value = 0.03398511728737064
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.882473
$ws.Range("N13").Value = 2.647419
$ws.Range("O13").Value = 0.1356390856021789
$ws.Range("P13").Value = 0.1356390856021789
$ws.Range("Q13").Value = 18.005247521883
$ws.Range("R13").Value = 162.047227696947
$ws.Range("S13").Value = 0.004609710232941755
$ws.Range("T13").Value = 0.004609710232941754
